$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update results to use consistent baseline for approach 3b
# Rows 6 (C naive3b weekly), 11 (P naive3b weekly), 16 (C naive3b monthly), 21 (P naive3b monthly)


# Row 6
$ws.Range("J6").Value = 18
$ws.Range("K6").Value = 24
$ws.Range("M6").Value = 30.75
$ws.Range("N6").Value = 36
$ws.Range("O6").Value = 27.23913043478261
$ws.Range("P6").Value = 4.571334755766127
$ws.Range("AL6").Value = 0.9473684210526315
$ws.Range("AN6").Value = 1.421052631578947
$ws.Range("AO6").Value = 1.631578947368421
$ws.Range("AP6").Value = 1.894736842105263
$ws.Range("AQ6").Value = 1.443036940176528
$ws.Range("AR6").Value = 0.240131076488877
$ws.Range("AS6").Value = 11.56756756756757
$ws.Range("AT6").Value = 13.375
$ws.Range("AU6").Value = 15.28571428571429
$ws.Range("AV6").Value = 17.12
$ws.Range("AW6").Value = 22.52631578947368
$ws.Range("AX6").Value = 15.51941047185875
$ws.Range("AY6").Value = 2.650506708276222

# Row 11
$ws.Range("K11").Value = 12
$ws.Range("L11").Value = 17
$ws.Range("N11").Value = 34
$ws.Range("O11").Value = 17.1304347826087
$ws.Range("P11").Value = 6.344231660517989
$ws.Range("AM11").Value = 0.631578947368421
$ws.Range("AN11").Value = 0.8947368421052632
$ws.Range("AP11").Value = 1.789473684210526
$ws.Range("AQ11").Value = 0.9073226544622426
$ws.Range("AR11").Value = 0.3332751163582248
$ws.Range("AS11").Value = 12.22857142857143
$ws.Range("AU11").Value = 23.77777777777778
$ws.Range("AV11").Value = 32.92307692307692
$ws.Range("AX11").Value = 26.84480643108547
$ws.Range("AY11").Value = 10.65341514253792

# Row 16
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 4
$ws.Range("M16").Value = 5
$ws.Range("N16").Value = 7
$ws.Range("O16").Value = 3.925925925925926
$ws.Range("P16").Value = 1.350366215671038
$ws.Range("S16").Value = 2
$ws.Range("T16").Value = 4
$ws.Range("U16").Value = 44
$ws.Range("V16").Value = 3.111111111111111
$ws.Range("W16").Value = 5.245875231005908
$ws.Range("AB16").Value = 33
$ws.Range("AC16").Value = 1.422222222222222
$ws.Range("AD16").Value = 4.317280340852534
$ws.Range("AJ16").Value = 0.137037037037037
$ws.Range("AK16").Value = 0.3755027837854114
$ws.Range("AL16").Value = 0
$ws.Range("AM16").Value = 1.5
$ws.Range("AN16").Value = 2
$ws.Range("AO16").Value = 2.5
$ws.Range("AP16").Value = 5
$ws.Range("AQ16").Value = 2.035802469135803
$ws.Range("AR16").Value = 0.7307731544163961
$ws.Range("AS16").Value = 11.83333333333333
$ws.Range("AT16").Value = 15.83333333333333
$ws.Range("AU16").Value = 19
$ws.Range("AV16").Value = 21.83333333333333
$ws.Range("AW16").Value = 59
$ws.Range("AX16").Value = 19.04585537918871
$ws.Range("AY16").Value = 5.78614337239353

# Row 21
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 4
$ws.Range("M21").Value = 5.25
$ws.Range("N21").Value = 8
$ws.Range("O21").Value = 4.216346153846154
$ws.Range("P21").Value = 1.726833441496639
$ws.Range("T21").Value = 4
$ws.Range("U21").Value = 21
$ws.Range("V21").Value = 3.177884615384615
$ws.Range("W21").Value = 3.512859928250906
$ws.Range("AB21").Value = 26
$ws.Range("AC21").Value = 2.418269230769231
$ws.Range("AD21").Value = 4.622257942559145
$ws.Range("AJ21").Value = 0.2644230769230769
$ws.Range("AK21").Value = 0.4937127935922568
$ws.Range("AM21").Value = 1.5
$ws.Range("AN21").Value = 2
$ws.Range("AO21").Value = 2.5
$ws.Range("AP21").Value = 4
$ws.Range("AQ21").Value = 2.012820512820513
$ws.Range("AR21").Value = 0.8083011042240419
$ws.Range("AS21").Value = 10.55555555555556
$ws.Range("AT21").Value = 14.26428571428572
$ws.Range("AU21").Value = 17.33333333333333
$ws.Range("AV21").Value = 19.63333333333333
$ws.Range("AW21").Value = 95
$ws.Range("AX21").Value = 20.93875534188034
$ws.Range("AY21").Value = 13.58581228733469
